$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before F ("Dop") which shifts old F..K to G..L ---
$ws.Columns("F").Insert()

# --- Column widths: E and the new F changed/created, old G (now H) changed too.
#     The rest (A-D, G, J, K, L) keep their widths automatically via the shift. ---
$ws.Columns("E").ColumnWidth = 15.166666666666666
$ws.Columns("F").ColumnWidth = 7
$ws.Columns("H").ColumnWidth = 15.166666666666666

# --- Apply the time-ish number formats first, in the exact order the original
#     authoring produced them, so the generated cellXfs come out as
#     1: h:mm, 2: h:mm:ss, 3: [h]:mm:ss (matching the target styles.xml). ---
$ws.Range("H2").NumberFormat = "h:mm"
$ws.Range("H3").NumberFormat = "h:mm"
$ws.Range("H4").NumberFormat = "h:mm"
$ws.Range("G4").NumberFormat = "h:mm"
$ws.Range("H6").NumberFormat = "h:mm"
$ws.Range("G2").NumberFormat = "h:mm:ss"
$ws.Range("G3").NumberFormat = "h:mm:ss"
$ws.Range("I2").NumberFormat = "h:mm:ss"
$ws.Range("I3").NumberFormat = "h:mm:ss"
$ws.Range("J2").NumberFormat = "h:mm:ss"
$ws.Range("J3").NumberFormat = "h:mm:ss"
$ws.Range("K2").NumberFormat = "h:mm:ss"
$ws.Range("K3").NumberFormat = "h:mm:ss"
$ws.Range("K4").NumberFormat = "h:mm:ss"
$ws.Range("K6").NumberFormat = "[h]:mm:ss"

# --- Plain numeric cells (do not add shared strings, order-independent) ---
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 337
$ws.Range("D2").Value = 500
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 144

$ws.Range("C3").Value = 337
$ws.Range("D3").Value = 500
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 144

$ws.Range("C4").Value = 337
$ws.Range("D4").Value = 500
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 144

$ws.Range("A6").Value = 2
$ws.Range("C6").Value = 337
$ws.Range("D6").Value = 500
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 144

$ws.Range("C7").Value = 337
$ws.Range("D7").Value = 500
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 144

$ws.Range("C8").Value = 337
$ws.Range("D8").Value = 500
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 144

# --- Text / shared-string cells, entered in the precise sequence that
#     reproduces the original sharedStrings.xml append order. ---
$ws.Range("F1").Value = "Dop"
$ws.Range("B2").Value = "a"
$ws.Range("B3").Value = "b"
$ws.Range("B4").Value = "c"
$ws.Range("L2").Value = "2,67 GB"
$ws.Range("I2").Value = "01:00:30h"
$ws.Range("I3").Value = "00:58:44h"
$ws.Range("H2").Value = "02:11m"
$ws.Range("H3").Value = "01:58m"
$ws.Range("H4").Value = "01:50m"
$ws.Range("G2").Value = "02:32:46h"
$ws.Range("G3").Value = "02:19:45h"
$ws.Range("G4").Value = "02:18:20h"
$ws.Range("I4").Value = "00:59:28h"
$ws.Range("J2").Value = "02:29:50h"
$ws.Range("J3").Value = "02:17:47h"
$ws.Range("K2").Value = "01:29:55h"
$ws.Range("K3").Value = "01:20:15h"
$ws.Range("K4").Value = "01:12:51h"
$ws.Range("J4").Value = "02:16:28h"
$ws.Range("L3").Value = "2,67 GB"
$ws.Range("L4").Value = "2,67 GB"
$ws.Range("G6").Value = "00:46:38h"
$ws.Range("H6").Value = "01:55m"
$ws.Range("I6").Value = "00:18:31h"
$ws.Range("J6").Value = "00:44:36h"
$ws.Range("K6").Value = "00:27:52h"
$ws.Range("L6").Value = "2,7GB"
$ws.Range("B6").Value = "a"
$ws.Range("B7").Value = "b"
$ws.Range("B8").Value = "c"

# --- Selection / active cell ---
$ws.Range("H7").Select()

# --- Workbook window view (best-effort; some attributes are controlled by
#     the host window manager and are not reachable from the object model) ---
$excel.Windows.Item(1).Height = 22760
